$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch to manual calculation so that the cells we are about to convert
# to text do not force a recalculation (and error-out) of the existing
# "SUM(C+D)" formulas that reference them - matches the source edit,
# which left those cached formula results untouched.
$excel.Calculation = -4135  # xlCalculationManual

# Rename header columns (also updates the Table1 column names & shared strings)
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Chai-Verkäufe insgesamt (Einheiten)"

# Replace numeric values in D6/E6 with text time-like values
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "17:05"
